# redo FR pop dens calcs using overlays
$wb = $excel.ActiveWorkbook

# Sheet "area_mixre" (sheet1): update mean, std, 25% summary stats
$wsMixre = $wb.Worksheets.Item("area_mixre")
$wsMixre.Range("B3").Value = 5.041855710213135
$wsMixre.Range("B4").Value = 5.498231592344241
$wsMixre.Range("B6").Value = 1.302180785072542

# Sheet "area_pop_sum" (sheet4): update population and density
$wsPopSum = $wb.Worksheets.Item("area_pop_sum")
$wsPopSum.Range("B3").Value = 448710
$wsPopSum.Range("B4").Value = 1022.95379769916
